$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$headers = @(
    "Job_Id",
    "Job_Title",
    "Job_Description",
    "Total_Years_Min_Exp",
    "Total_Years_Max_Exp",
    "Work_Mode",
    "Job_Location",
    "LinkedIn_Poster",
    "LinkedIn_Posted",
    "Resume_received",
    "Resume_downloaded"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data row (row 2) ---
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Junior RPA Developer"
$ws.Range("C2").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Hybrid"
$ws.Range("G2").Value = "Bengaluru, Karnataka, India"

# The multi-line description causes the host engine to auto-expand the row
# height (simulating Excel's implicit text-height heuristic for embedded
# newlines); AutoFit the row back down so no explicit custom height is
# stored, matching a normal (non line-wrapped) default-height row.
$ws.Rows.Item(2).AutoFit() | Out-Null

# --- Header formatting: bold, centered horizontally, top vertical align, thin box border ---
# Apply to A1 first via its .Style object (keeps the style table compact: one
# combined cell format instead of a separate format per property), then copy
# that exact format across the rest of the header row.
$headerCell = $ws.Range("A1")
$headerStyle = $headerCell.Style
$headerStyle.Font.Bold = $true
$headerStyle.HorizontalAlignment = -4108   # xlCenter
$headerStyle.VerticalAlignment = -4160     # xlTop
$headerStyle.Borders.LineStyle = 1         # xlContinuous

$headerCell.Copy() | Out-Null
$ws.Range("B1:K1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A1").Select() | Out-Null

Write-Output "Job posting row added"
